$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.720.60"
$ws.Range("E2").Value = "  +3.12%  "
$ws.Range("D3").Value = "2.444.69"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "2.443.23"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  +7.74%  "
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "62.583.01"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.446.25"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("B19").Value = "BabyDogeCoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D19").Value = "0.0₆0956"
$ws.Range("E19").Value = "  +245.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "329.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +9.96%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "638.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.12%  "
$ws.Range("E28").Value = "  +16.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.11%  "
$ws.Range("D30").Value = "0.0₃0983"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").Value = "2.565.80"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  +7.83%  "
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("E35").Value = "  +4.45%  "
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.37%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.08%  "
$ws.Range("E44").Value = "  +4.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +27.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.605"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.93%  "
